$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the scraped cryptos-list values cell by cell (Coin / Link / Price / Volume(1h)).
# Price-column values that look like plain numbers ("237.46", "1.001", ...) need to be
# forced to Text before assignment so Excel keeps them as strings (matching the original
# inlineStr cells) instead of silently converting them to numeric values; the temporary
# Text number-format is cleared right back off again afterwards so the cell style is left
# untouched overall.

$ws.Range("D2").Value = '25.820.59'
$ws.Range("E2").Value = '  +0.54%  '
$ws.Range("D3").Value = '1.757.71'
$ws.Range("E3").Value = '  +0.82%  '
$ws.Range("E4").Value = '  -0.14%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '237.46'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -0.10%  '
$ws.Range("E6").Value = '  -0.25%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5076'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +3.31%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '41.20'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -0.02%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2643'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +8.78%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.06208'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +3.40%  '
$ws.Range("D11").Value = '1.753.32'
$ws.Range("E11").Value = '  +0.47%  '
$ws.Range("E12").Value = '  +4.56%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '15.56'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +8.84%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6044'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +1.59%  '
$ws.Range("E15").Value = '  +3.17%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '77.52'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +0.83%  '
$ws.Range("E17").Value = '  -0.13%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.001'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -0.30%  '
$ws.Range("D19").Value = '25.861.72'
$ws.Range("E19").Value = '  +0.56%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000006838'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +9.25%  '
$ws.Range("E21").Value = '  +4.99%  '
$ws.Range("D22").Value = '1.978.79'
$ws.Range("E22").Value = '  -0.52%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.068'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +5.99%  '
$ws.Range("B24").Value = 'Cosmos'
$ws.Range("C24").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '8.158'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +2.76%  '
$ws.Range("B25").Value = 'Chainlink'
$ws.Range("C25").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.178'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +2.07%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '137.84'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +2.73%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.455'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +2.21%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.821'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -2.32%  '
$ws.Range("E29").Value = '  +5.47%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '102.63'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +3.07%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08237'
$ws.Range("D31").ClearFormats()
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.699'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +2.74%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.401'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +7.87%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.04364'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +2.50%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.000'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -0.32%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.656'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +1.63%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.002'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -1.62%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.6005'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -1.08%  '
$ws.Range("E39").Value = '  +1.08%  '
$ws.Range("E40").Value = '  +6.58%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.933'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -7.80%  '
$ws.Range("E42").Value = '  -0.24%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '103.37'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +2.06%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.3838'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +1.09%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.7444'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -5.22%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '4.876'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -4.92%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.05493'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +8.44%  '
$ws.Range("E48").Value = '  +4.90%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '5.960'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -2.22%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '30.12'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +3.12%  '
$ws.Range("B51").Value = 'USDD'
$ws.Range("C51").Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.001'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +0.19%  '
